# Update "想去人数" (want-to-go count) values in column F across sheets,
# matching the scraped data refresh described in the commit.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 292
$ws1.Range("F5").Value  = 770
$ws1.Range("F6").Value  = 464
$ws1.Range("F11").Value = 7176
$ws1.Range("F14").Value = 1402
$ws1.Range("F17").Value = 399
$ws1.Range("F21").Value = 732
$ws1.Range("F23").Value = 46
$ws1.Range("F24").Value = 120
$ws1.Range("F26").Value = 202
$ws1.Range("F30").Value = 1059
$ws1.Range("F32").Value = 81
$ws1.Range("F33").Value = 2056
$ws1.Range("F34").Value = 577
$ws1.Range("F36").Value = 15
$ws1.Range("F38").Value = 556

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value  = 281
$ws2.Range("F4").Value  = 60
$ws2.Range("F10").Value = 137

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 292
$ws4.Range("F6").Value  = 770
$ws4.Range("F8").Value  = 464
$ws4.Range("F13").Value = 7176
$ws4.Range("F16").Value = 281
$ws4.Range("F17").Value = 1402
$ws4.Range("F20").Value = 399
$ws4.Range("F22").Value = 60
$ws4.Range("F28").Value = 732
$ws4.Range("F30").Value = 46
$ws4.Range("F31").Value = 120
$ws4.Range("F35").Value = 137
$ws4.Range("F36").Value = 202
$ws4.Range("F40").Value = 1059
$ws4.Range("F42").Value = 81
$ws4.Range("F43").Value = 2056
$ws4.Range("F44").Value = 577
$ws4.Range("F46").Value = 15
$ws4.Range("F48").Value = 556
